$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Nxph3"
$ws.Range("C2").Value = "Nrxn1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1405946666666667
$ws.Range("H2").Value = 0.421784
$ws.Range("I2").Value = 0.539087912495223
$ws.Range("J2").Value = 0.539087912495223
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.006517333333333333
$ws.Range("N2").Value = 0.019552
$ws.Range("O2").Value = 0.004697037253424763
$ws.Range("P2").Value = 0.004697037253424762
$ws.Range("Q2").Value = 0.0009163023075555556
$ws.Range("R2").Value = 0.008246720768
$ws.Range("S2").Value = 0.002532116007861051
$ws.Range("T2").Value = 0.002532116007861051

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Nxph3"
$ws.Range("C3").Value = "Nrxn1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1405946666666667
$ws.Range("H3").Value = 0.421784
$ws.Range("I3").Value = 0.539087912495223
$ws.Range("J3").Value = 0.539087912495223
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01189366666666667
$ws.Range("N3").Value = 0.035681
$ws.Range("O3").Value = 0.008571756661182945
$ws.Range("P3").Value = 0.008571756661182945
$ws.Range("Q3").Value = 0.001672186100444444
$ws.Range("R3").Value = 0.015049674904
$ws.Range("S3").Value = 0.004620930404894137
$ws.Range("T3").Value = 0.004620930404894137

# Row 4
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Nxph3"
$ws.Range("C4").Value = "Nrxn1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1405946666666667
$ws.Range("H4").Value = 0.421784
$ws.Range("I4").Value = 0.539087912495223
$ws.Range("J4").Value = 0.539087912495223
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.369130333333333
$ws.Range("N4").Value = 4.107391
$ws.Range("O4").Value = 0.9867312060853923
$ws.Range("P4").Value = 0.9867312060853922
$ws.Range("Q4").Value = 0.1924924228382222
$ws.Range("R4").Value = 1.732431805544
$ws.Range("S4").Value = 0.5319348660824679
$ws.Range("T4").Value = 0.5319348660824678

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Nxph3"
$ws.Range("C5").Value = "Nrxn1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1202063333333333
$ws.Range("H5").Value = 0.360619
$ws.Range("I5").Value = 0.460912087504777
$ws.Range("J5").Value = 0.460912087504777
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.006517333333333333
$ws.Range("N5").Value = 0.019552
$ws.Range("O5").Value = 0.004697037253424763
$ws.Range("P5").Value = 0.004697037253424762
$ws.Range("Q5").Value = 0.0007834247431111112
$ws.Range("R5").Value = 0.007050822688000001
$ws.Range("S5").Value = 0.002164921245563712
$ws.Range("T5").Value = 0.002164921245563711

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Nxph3"
$ws.Range("C6").Value = "Nrxn1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1202063333333333
$ws.Range("H6").Value = 0.360619
$ws.Range("I6").Value = 0.460912087504777
$ws.Range("J6").Value = 0.460912087504777
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.01189366666666667
$ws.Range("N6").Value = 0.035681
$ws.Range("O6").Value = 0.008571756661182945
$ws.Range("P6").Value = 0.008571756661182945
$ws.Range("Q6").Value = 0.001429694059888889
$ws.Range("R6").Value = 0.012867246539
$ws.Range("S6").Value = 0.003950826256288808
$ws.Range("T6").Value = 0.003950826256288809

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Nxph3"
$ws.Range("C7").Value = "Nrxn1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.1202063333333333
$ws.Range("H7").Value = 0.360619
$ws.Range("I7").Value = 0.460912087504777
$ws.Range("J7").Value = 0.460912087504777
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.369130333333333
$ws.Range("N7").Value = 4.107391
$ws.Range("O7").Value = 0.9867312060853923
$ws.Range("P7").Value = 0.9867312060853922
$ws.Range("Q7").Value = 0.1645781372254445
$ws.Range("R7").Value = 1.481203235029
$ws.Range("S7").Value = 0.4547963400029245
$ws.Range("T7").Value = 0.4547963400029245
